$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current values from row 2 (columns D, I, J, K, L, M, P)
$d2 = $ws.Range("D2").Value2
$i2 = $ws.Range("I2").Value2
$j2 = $ws.Range("J2").Value2
$k2 = $ws.Range("K2").Value2
$l2 = $ws.Range("L2").Value2
$m2 = $ws.Range("M2").Value2
$p2 = $ws.Range("P2").Value2

# Read current values from row 5 (columns D, I, J, K, L, M, P)
$d5 = $ws.Range("D5").Value2
$i5 = $ws.Range("I5").Value2
$j5 = $ws.Range("J5").Value2
$k5 = $ws.Range("K5").Value2
$l5 = $ws.Range("L5").Value2
$m5 = $ws.Range("M5").Value2
$p5 = $ws.Range("P5").Value2

# Swap row 2 values with row 5 values
$ws.Range("D2").Value2 = $d5
$ws.Range("I2").Value2 = $i5
$ws.Range("J2").Value2 = $j5
$ws.Range("K2").Value2 = $k5
$ws.Range("L2").Value2 = $l5
$ws.Range("M2").Value2 = $m5
$ws.Range("P2").Value2 = $p5

$ws.Range("D5").Value2 = $d2
$ws.Range("I5").Value2 = $i2
$ws.Range("J5").Value2 = $j2
$ws.Range("K5").Value2 = $k2
$ws.Range("L5").Value2 = $l2
$ws.Range("M5").Value2 = $m2
$ws.Range("P5").Value2 = $p2
